$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value2 = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value2 = '27.415.44'
$ws.Range('E2').Value2 = '  +2.08%  '
$ws.Range('D3').Value2 = '1.828.09'
$ws.Range('E3').Value2 = '  +1.07%  '
Set-TextValue 'D4' '1.0000'
$ws.Range('E4').Value2 = '  -0.09%  '
Set-TextValue 'D5' '313.54'
$ws.Range('E5').Value2 = '  +1.16%  '
Set-TextValue 'D6' '1.0000'
Set-TextValue 'D7' '0.4451'
$ws.Range('E7').Value2 = '  -0.02%  '
Set-TextValue 'D8' '0.3775'
$ws.Range('E8').Value2 = '  +2.88%  '
Set-TextValue 'D9' '0.07406'
$ws.Range('E9').Value2 = '  +1.74%  '
Set-TextValue 'D10' '0.8798'
$ws.Range('E10').Value2 = '  +3.28%  '
$ws.Range('E11').Value2 = '  +1.30%  '
$ws.Range('D12').Value2 = '1.833.54'
$ws.Range('E12').Value2 = '  +1.47%  '
Set-TextValue 'D13' '6.725'
$ws.Range('E13').Value2 = '  +2.02%  '
Set-TextValue 'D14' '5.432'
$ws.Range('E14').Value2 = '  +2.63%  '
Set-TextValue 'D15' '92.97'
$ws.Range('E15').Value2 = '  +1.62%  '
$ws.Range('E16').Value2 = '  -0.14%  '
$ws.Range('E17').Value2 = '  -0.16%  '
Set-TextValue 'D18' '0.000008809'
$ws.Range('E18').Value2 = '  +1.12%  '
Set-TextValue 'D19' '1.0000'
$ws.Range('E19').Value2 = '  -0.07%  '
$ws.Range('E20').Value2 = '  +1.40%  '
$ws.Range('D21').Value2 = '27.425.03'
$ws.Range('E21').Value2 = '  +2.07%  '
Set-TextValue 'D22' '5.358'
$ws.Range('E22').Value2 = '  +4.25%  '
Set-TextValue 'D23' '10.97'
$ws.Range('E23').Value2 = '  +1.48%  '
Set-TextValue 'D24' '1.947'
$ws.Range('E24').Value2 = '  -1.64%  '
Set-TextValue 'D25' '151.09'
$ws.Range('E25').Value2 = '  -0.07%  '
Set-TextValue 'D26' '2.286'
$ws.Range('E26').Value2 = '  +3.23%  '
Set-TextValue 'D27' '18.69'
$ws.Range('E27').Value2 = '  +1.87%  '
Set-TextValue 'D28' '5.362'
$ws.Range('E28').Value2 = '  +3.26%  '
Set-TextValue 'D29' '117.26'
$ws.Range('E29').Value2 = '  +0.93%  '
Set-TextValue 'D30' '0.08909'
$ws.Range('E30').Value2 = '  +1.29%  '
Set-TextValue 'D31' '0.7935'
$ws.Range('E31').Value2 = '  +6.30%  '
Set-TextValue 'D32' '1.201'
$ws.Range('E32').Value2 = '  +2.46%  '
Set-TextValue 'D33' '4.565'
$ws.Range('E33').Value2 = '  +3.04%  '
Set-TextValue 'D34' '2.939'
$ws.Range('E34').Value2 = '  +0.20%  '
Set-TextValue 'D35' '0.9996'
$ws.Range('E35').Value2 = '  -0.06%  '
Set-TextValue 'D36' '1.109'
$ws.Range('E36').Value2 = '  +1.57%  '
Set-TextValue 'D37' '0.01985'
$ws.Range('E37').Value2 = '  +1.42%  '
Set-TextValue 'D38' '0.05274'
$ws.Range('E38').Value2 = '  +1.80%  '
Set-TextValue 'D39' '7.303'
$ws.Range('E39').Value2 = '  +3.62%  '
Set-TextValue 'D40' '0.5328'
$ws.Range('B41').Value2 = 'RenderToken'
$ws.Range('C41').Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D41' '2.354'
$ws.Range('E41').Value2 = '  +19.82%  '
$ws.Range('B42').Value2 = 'MXToken'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D42' '2.875'
$ws.Range('E42').Value2 = '  +0.07%  '
Set-TextValue 'D43' '0.1706'
$ws.Range('E43').Value2 = '  +1.27%  '
Set-TextValue 'D44' '8.712'
$ws.Range('E44').Value2 = '  +3.48%  '
Set-TextValue 'D45' '0.5067'
$ws.Range('E45').Value2 = '  -1.83%  '
Set-TextValue 'D46' '10.62'
$ws.Range('E46').Value2 = '  +1.28%  '
Set-TextValue 'D47' '105.60'
$ws.Range('E47').Value2 = '  +0.22%  '
$ws.Range('E48').Value2 = '  +2.34%  '
Set-TextValue 'D49' '0.9995'
$ws.Range('E49').Value2 = '  -0.04%  '
Set-TextValue 'D50' '0.06392'
$ws.Range('E50').Value2 = '  +1.08%  '
Set-TextValue 'D51' '66.34'
$ws.Range('E51').Value2 = '  +6.17%  '
